$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 42612.883032407408
$ws.Range("B10").Value = -8
$ws.Range("C10").Value = 47
$ws.Range("D10").Value = 51
$ws.Range("E10").Value = 40
$ws.Range("F10").Value = 60
$ws.Range("G10").Value = 14722
$ws.Range("H10").Value = 13875
$ws.Range("I10").Value = 881
$ws.Range("J10").Value = 129
$ws.Range("K10").Value = 138
$ws.Range("L10").Value = 2
$ws.Range("M10").Value = 3
$ws.Range("N10").Value = "Named"

$ws.Range("A11").Value = 42612.88958333333
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = 48
$ws.Range("D11").Value = 49
$ws.Range("E11").Value = 65
$ws.Range("F11").Value = 34
$ws.Range("G11").Value = 19463
$ws.Range("H11").Value = 15316
$ws.Range("I11").Value = 973
$ws.Range("J11").Value = 156
$ws.Range("K11").Value = 160
$ws.Range("L11").Value = 15
$ws.Range("M11").Value = 8
$ws.Range("N11").Value = "Named"
